# Scheduled-runner market data refresh: update currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H-N) on affected leve rows across
# the Sheets workbook, as produced by the latest data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 513.7765000000001
$ws.Range("I15").Value = 513.7765000000001
$ws.Range("K15").Value = 1541.3295
$ws.Range("M15").Value = -1372.3295
$ws.Range("H17").Value = 2615576.5
$ws.Range("J17").Value = 2615576.5
$ws.Range("L17").Value = 7846729.5
$ws.Range("N17").Value = -7847065.5
$ws.Range("H33").Value = 9262491
$ws.Range("I33").Value = 15625819
$ws.Range("K33").Value = 15625819
$ws.Range("M33").Value = -15625590
$ws.Range("H41").Value = 197.58333
$ws.Range("I41").Value = 127.1
$ws.Range("K41").Value = 127.1
$ws.Range("M41").Value = 312.9
$ws.Range("H53").Value = 16890.285
$ws.Range("J53").Value = 28033
$ws.Range("L53").Value = 28033
$ws.Range("N53").Value = -29307
$ws.Range("H55").Value = 190.77777
$ws.Range("J55").Value = 177.33333
$ws.Range("L55").Value = 177.33333
$ws.Range("N55").Value = -605.3333299999999
$ws.Range("H112").Value = 2606.68
$ws.Range("J112").Value = 2827.9524
$ws.Range("L112").Value = 8483.8572
$ws.Range("N112").Value = -10699.8572
$ws.Range("H123").Value = 167599.8
$ws.Range("J123").Value = 167599.8
$ws.Range("L123").Value = 167599.8
$ws.Range("N123").Value = -177399.8
$ws.Range("H125").Value = 9818.909
$ws.Range("I125").Value = 6418.1665
$ws.Range("K125").Value = 57763.4985
$ws.Range("M125").Value = -55303.4985
$ws.Range("H132").Value = 10064.648
$ws.Range("I132").Value = 7837.4375
$ws.Range("J132").Value = 24318.8
$ws.Range("K132").Value = 23512.3125
$ws.Range("L132").Value = 72956.39999999999
$ws.Range("M132").Value = -20982.3125
$ws.Range("N132").Value = -78016.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1537301
$ws.Range("I32").Value = 4357599
$ws.Range("K32").Value = 4357599
$ws.Range("M32").Value = -4357312
$ws.Range("H61").Value = 11419.345
$ws.Range("J61").Value = 31325.223
$ws.Range("L61").Value = 31325.223
$ws.Range("N61").Value = -31749.223
$ws.Range("H74").Value = 24375.578
$ws.Range("I74").Value = 3227
$ws.Range("J74").Value = 30015.2
$ws.Range("K74").Value = 3227
$ws.Range("L74").Value = 30015.2
$ws.Range("M74").Value = -2353
$ws.Range("N74").Value = -31763.2
$ws.Range("H77").Value = 24375.578
$ws.Range("I77").Value = 3227
$ws.Range("J77").Value = 30015.2
$ws.Range("K77").Value = 16135
$ws.Range("L77").Value = 150076
$ws.Range("M77").Value = -11767
$ws.Range("N77").Value = -158812
$ws.Range("H97").Value = 2786.074
$ws.Range("I97").Value = 900.6957
$ws.Range("J97").Value = 13627
$ws.Range("K97").Value = 900.6957
$ws.Range("L97").Value = 13627
$ws.Range("M97").Value = -404.6957
$ws.Range("N97").Value = -14619
$ws.Range("H136").Value = 11419.345
$ws.Range("J136").Value = 31325.223
$ws.Range("L136").Value = 93975.66900000001
$ws.Range("N136").Value = -99075.66900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 25340.967
$ws.Range("I20").Value = 6159
$ws.Range("J20").Value = 45801.734
$ws.Range("K20").Value = 6159
$ws.Range("L20").Value = 45801.734
$ws.Range("M20").Value = -5912
$ws.Range("N20").Value = -46295.734
$ws.Range("H134").Value = 8146.636
$ws.Range("I134").Value = 2578
$ws.Range("J134").Value = 15473.789
$ws.Range("K134").Value = 7734
$ws.Range("L134").Value = 46421.367
$ws.Range("M134").Value = -5199
$ws.Range("N134").Value = -51491.367

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 368.08334
$ws.Range("I7").Value = 119.117645
$ws.Range("K7").Value = 119.117645
$ws.Range("M7").Value = -6.117644999999996
$ws.Range("H33").Value = 5498.3335
$ws.Range("I33").Value = 5498.3335
$ws.Range("K33").Value = 5498.3335
$ws.Range("M33").Value = -5119.3335
$ws.Range("H86").Value = 4699
$ws.Range("I86").Value = 3866.2727
$ws.Range("J86").Value = 5531.727
$ws.Range("K86").Value = 3866.2727
$ws.Range("L86").Value = 5531.727
$ws.Range("M86").Value = -2743.2727
$ws.Range("N86").Value = -7777.727
$ws.Range("H89").Value = 4699
$ws.Range("I89").Value = 3866.2727
$ws.Range("J89").Value = 5531.727
$ws.Range("K89").Value = 19331.3635
$ws.Range("L89").Value = 27658.635
$ws.Range("M89").Value = -13715.3635
$ws.Range("N89").Value = -38890.63499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 3050.5
$ws.Range("I11").Value = 3050.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 9151.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -9011.5
$ws.Range("N11").ClearContents()
$ws.Range("H12").Value = 29.555555
$ws.Range("J12").Value = 50.5
$ws.Range("L12").Value = 151.5
$ws.Range("N12").Value = -497.5
$ws.Range("H40").Value = 548.8461
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H68").Value = 1163.5454
$ws.Range("J68").Value = 2019.8
$ws.Range("L68").Value = 6059.4
$ws.Range("N68").Value = -7681.4
$ws.Range("H71").Value = 1163.5454
$ws.Range("J71").Value = 2019.8
$ws.Range("L71").Value = 18178.2
$ws.Range("N71").Value = -26290.2
$ws.Range("H75").Value = 4176.6665
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 5765
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 17295
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -19291
$ws.Range("H78").Value = 4176.6665
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 5765
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 51885
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -61869
$ws.Range("H92").Value = 1491.3077
$ws.Range("I92").Value = 1799.6666
$ws.Range("J92").Value = 1398.8
$ws.Range("K92").Value = 5398.9998
$ws.Range("L92").Value = 4196.4
$ws.Range("M92").Value = -4150.9998
$ws.Range("N92").Value = -6692.4
$ws.Range("H129").Value = 8334633
$ws.Range("J129").Value = 1462.125
$ws.Range("L129").Value = 4386.375
$ws.Range("N129").Value = -14386.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1036.238
$ws.Range("I2").Value = 1249.5294
$ws.Range("K2").Value = 1249.5294
$ws.Range("M2").Value = -1136.5294
$ws.Range("H70").Value = 8126.654
$ws.Range("I70").Value = 5985.6924
$ws.Range("J70").Value = 10267.615
$ws.Range("K70").Value = 5985.6924
$ws.Range("L70").Value = 10267.615
$ws.Range("M70").Value = -5715.6924
$ws.Range("N70").Value = -10807.615
$ws.Range("H73").Value = 8126.654
$ws.Range("I73").Value = 5985.6924
$ws.Range("J73").Value = 10267.615
$ws.Range("K73").Value = 5985.6924
$ws.Range("L73").Value = 10267.615
$ws.Range("M73").Value = -5049.6924
$ws.Range("N73").Value = -12139.615
$ws.Range("H107").Value = 1933
$ws.Range("J107").Value = 1900
$ws.Range("L107").Value = 1900
$ws.Range("N107").Value = -5740
$ws.Range("H113").Value = 92044.39999999999
$ws.Range("J113").Value = 2998.2
$ws.Range("L113").Value = 2998.2
$ws.Range("N113").Value = -7338.2
$ws.Range("H122").Value = 6225.1763
$ws.Range("I122").Value = 2777.1667
$ws.Range("K122").Value = 8331.500100000001
$ws.Range("M122").Value = -5881.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11752.77
$ws.Range("I40").Value = 7598.375
$ws.Range("J40").Value = 18399.8
$ws.Range("K40").Value = 7598.375
$ws.Range("L40").Value = 18399.8
$ws.Range("M40").Value = -7462.375
$ws.Range("N40").Value = -18671.8
$ws.Range("H46").Value = 2002879.8
$ws.Range("I46").Value = 5000450.5
$ws.Range("K46").Value = 5000450.5
$ws.Range("M46").Value = -5000262.5
$ws.Range("H93").Value = 6394.6343
$ws.Range("I93").Value = 3804.5483
$ws.Range("J93").Value = 14423.9
$ws.Range("K93").Value = 3804.5483
$ws.Range("L93").Value = 14423.9
$ws.Range("M93").Value = -2556.5483
$ws.Range("N93").Value = -16919.9
$ws.Range("H122").Value = 7558.12
$ws.Range("I122").Value = 5113.6665
$ws.Range("J122").Value = 11224.8
$ws.Range("K122").Value = 15340.9995
$ws.Range("L122").Value = 33674.39999999999
$ws.Range("M122").Value = -12890.9995
$ws.Range("N122").Value = -38574.39999999999
$ws.Range("H136").Value = 11866.891
$ws.Range("I136").Value = 11385.167
$ws.Range("J136").Value = 12444.96
$ws.Range("K136").Value = 34155.501
$ws.Range("L136").Value = 37334.88
$ws.Range("M136").Value = -31605.501
$ws.Range("N136").Value = -42434.88

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3952.35
$ws.Range("J113").Value = 7808.5557
$ws.Range("L113").Value = 23425.6671
$ws.Range("N113").Value = -27765.6671
